# Correcciones y modificaciones en sprint backlog 1, aun queda pendiente sprint backlog 2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT BACKLOG 1")
$ws.Activate() | Out-Null

# --- Fill in the daily burn-down numbers for the Sprint Backlog 1 tasks (rows 6-10) ---
# Row 6: "Desarrollo del login por lado cliente"
$ws.Cells.Item(6, 4).Value = 5   # D6
$ws.Cells.Item(6, 5).Value = 2   # E6
$ws.Cells.Item(6, 6).Value = 0   # F6
$ws.Cells.Item(6, 7).Value = 0   # G6
$ws.Cells.Item(6, 8).Value = 0   # H6
$ws.Cells.Item(6, 9).Value = 0   # I6

# Row 7: "Desarrollo del login por lado servidor"
$ws.Cells.Item(7, 4).Value = 8   # D7
$ws.Cells.Item(7, 5).Value = 8   # E7
$ws.Cells.Item(7, 6).Value = 6   # F7
$ws.Cells.Item(7, 7).Value = 4   # G7
$ws.Cells.Item(7, 8).Value = 2   # H7
$ws.Cells.Item(7, 9).Value = 0   # I7

# Row 8: "Diseño Login UI de la aplicación" -- D8:H8 were already 0, only I8 was blank
$ws.Cells.Item(8, 9).Value = 0   # I8

# Row 9: "Pruebas funcionales por lado del cliente"
$ws.Cells.Item(9, 4).Value = 8   # D9
$ws.Cells.Item(9, 5).Value = 7   # E9
$ws.Cells.Item(9, 6).Value = 3   # F9
$ws.Cells.Item(9, 7).Value = 3   # G9
$ws.Cells.Item(9, 8).Value = 0   # H9
$ws.Cells.Item(9, 9).Value = 0   # I9

# Row 10: "Pruebas por lado del servidor" -- D10:H10 were already 8, only I10 was blank
$ws.Cells.Item(10, 9).Value = 0  # I10

# --- Update the authorship note: it used to be a single "Autor: ..." line,
#     now it is a bold "Autores:" heading followed by one line per author ---
$ws.Range("A22").Value = "Autores:"
$ws.Range("A22").Font.Bold = $true

$ws.Range("A23").Value = "Luis Gianpierre Portella Bravo"
$ws.Range("A24").Value = "Hans Soto Rojas"

# --- Update the view so the newly added rows are visible/selected like in the
#     authored workbook ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select() | Out-Null

Write-Output "Sprint Backlog 1 updated: burn-down numbers filled in and author list expanded."
